# Adds a new sale-record row for "رول اون ريكسونا25" right before the
# existing "سرنجات 3 سم" row (keeping the A-column numbering + totals
# consistent), and refreshes the generated-at timestamp in the footer.
#
# Strategy: work from the bottom of the affected block upward so that each
# row's original content is read/copied before it gets overwritten.
#
#   row 58 (footer)      -> copied down to row 59
#   row 57 (totals)      -> copied down to row 58, total bumped by 35
#   row 56 (سرنجات 3 سم) -> copied down to row 57 (unchanged values)
#   row 56               -> overwritten in place with the new product

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: footer row 58 -> row 59
# ---------------------------------------------------------------------
$ws.Range("A58:Q58").Copy()
$ws.Range("A59:Q59").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A58:F58").UnMerge()
$ws.Range("G58:I58").UnMerge()
$ws.Range("K58:Q58").UnMerge()

$ws.Cells.Item(59, 1).Value = "Friday, 25 July, 2025 9:13 PM"
$ws.Cells.Item(59, 7).Value = "1/1"
$ws.Cells.Item(59, 11).Value = "developed by : Abdelaziz Talaat"

$ws.Range("A59:F59").Merge()
$ws.Range("G59:I59").Merge()
$ws.Range("K59:Q59").Merge()
$ws.Rows(59).RowHeight = 16.5

# ---------------------------------------------------------------------
# Step 2: totals row 57 -> row 58 (grand total grows by the new line's
# sale price, 35.0000)
# ---------------------------------------------------------------------
$ws.Range("P57:Q57").Copy()
$ws.Range("P58:Q58").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("P57:Q57").UnMerge()

$ws.Cells.Item(58, 16).Value = 3181.9749999999999
$ws.Cells.Item(58, 17).Value = ""

$ws.Range("P58:Q58").Merge()
$ws.Rows(58).RowHeight = 24.75

# ---------------------------------------------------------------------
# Step 3: product row 56 -> row 57 (سرنجات 3 سم keeps its old values,
# just moves one row down and its serial number becomes 51)
# ---------------------------------------------------------------------
$ws.Range("A56:Q56").Copy()
$ws.Range("A57:Q57").PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(57, 1).Value = 51
$ws.Cells.Item(57, 3).Value = "سرنجات 3 سم"
$ws.Cells.Item(57, 8).Value = "0:0"
$ws.Cells.Item(57, 12).Value = "0"
$ws.Cells.Item(57, 14).Value = "2.00"
$ws.Cells.Item(57, 16).Value = "10.0000"
$ws.Cells.Item(57, 17).Value = "5:0"

$ws.Range("A57:B57").Merge()
$ws.Range("C57:G57").Merge()
$ws.Range("H57:K57").Merge()
$ws.Range("L57:M57").Merge()
$ws.Range("N57:O57").Merge()
$ws.Rows(57).RowHeight = 25.5

# ---------------------------------------------------------------------
# Step 4: row 56 becomes the newly-added product, رول اون ريكسونا25
# (serial number 50 and the L-column "حد الطلب" value of 0 stay as-is)
# ---------------------------------------------------------------------
$ws.Cells.Item(56, 3).Value = "رول اون ريكسونا25"
$ws.Cells.Item(56, 8).Value = "4:0"
$ws.Cells.Item(56, 14).Value = "35.00"
$ws.Cells.Item(56, 16).Value = "35.0000"
$ws.Cells.Item(56, 17).Value = "1:0"

Write-Output "edit applied"
